# Weekly update: insert the new "Brócoli" price record for
# Vega Monumental Concepción at the top of the dated block (row 496),
# shifting the existing rows 496:586 down to 497:587.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 496 (pushes old 496..586 down to 497..587,
# extends the used range to A1:R587, same as the target dimension).
$ws.Rows.Item(496).Insert()

# Populate the new row with this week's record.
$ws.Range("A496").Value = 11
$ws.Range("B496").Value = "Vega Monumental Concepción"
$ws.Range("C496").Value = "Bíobío"
$ws.Range("D496").Value = 45209
$ws.Range("E496").Value = 8
$ws.Range("F496").Value = 100112023
$ws.Range("G496").Value = "Brócoli"
$ws.Range("H496").Value = "Sin especificar"
$ws.Range("I496").Value = "Primera"
$ws.Range("J496").Value = 2000
$ws.Range("K496").Value = 800
$ws.Range("L496").Value = 900
$ws.Range("M496").Value = 850
$ws.Range("N496").Value = "$/unidad"
$ws.Range("O496").Value = "Región Metropolitana"
$ws.Range("P496").Value = 850
$ws.Range("Q496").Value = 1
$ws.Range("R496").Value = "Hortaliza"
